# ajustes finais para a nova versao
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2021-04-28"
$ws.Range("D2").Value = "entrada"
$ws.Range("E2").Value = 21.21

# Row 3
$ws.Range("A3").Value = "Transporte"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2021-04-27"
$ws.Range("D3").Value = "saida"
$ws.Range("E3").Value = 178.8

# Row 4 (Total Geral)
$ws.Range("E4").Value = -157.59
